# Auto-update stock values: 2025-12-13 07:53:11 UTC
#
# Appends a new trading-day column (date 20251212) to every data sheet of
# the workbook, one column to the right of the current last column,
# mirroring the header/number formatting of the existing last column and
# filling in the three data rows (row1 = date header, row2/row3 = values).

$wb = $excel.ActiveWorkbook

# xlPasteFormats / xlPasteValues paste-special constants used below.
$xlPasteFormats = -4122
$xlPasteValues  = -4163

# One entry per worksheet that grows by a column in this commit. OldCol /
# NewCol are 1-based column indexes (OldCol = current last column,
# NewCol = OldCol + 1, the column being appended).  Row1Type is "n" for
# sheets that store the date header as a genuine number and "s" for the
# sheets that store it as text (matches the source workbook exactly).
$entries = @(
    @{ SheetIdx = 2;  Name = "시가";   OldCol = 75; NewCol = 76; Row1Type = "n"; Row1Val = "20251212"; Row2Val = 622.08;    Row3Val = 55.11 },
    @{ SheetIdx = 3;  Name = "고가";   OldCol = 75; NewCol = 76; Row1Type = "n"; Row1Val = "20251212"; Row2Val = 623.54;    Row3Val = 55.52 },
    @{ SheetIdx = 4;  Name = "저가";   OldCol = 75; NewCol = 76; Row1Type = "n"; Row1Val = "20251212"; Row2Val = 611.36;    Row3Val = 52.23 },
    @{ SheetIdx = 5;  Name = "종가";   OldCol = 75; NewCol = 76; Row1Type = "n"; Row1Val = "20251212"; Row2Val = 613.62;    Row3Val = 52.82 },
    @{ SheetIdx = 6;  Name = "거래량"; OldCol = 75; NewCol = 76; Row1Type = "n"; Row1Val = "20251212"; Row2Val = 75158726; Row3Val = 138448531 },
    @{ SheetIdx = 7;  Name = "s20";   OldCol = 56; NewCol = 57; Row1Type = "n"; Row1Val = "20251212"; Row2Val = 67;        Row3Val = 11 },
    @{ SheetIdx = 8;  Name = "s60";   OldCol = 16; NewCol = 17; Row1Type = "n"; Row1Val = "20251212"; Row2Val = 56;        Row3Val = 9 },
    @{ SheetIdx = 9;  Name = "z20";   OldCol = 56; NewCol = 57; Row1Type = "n"; Row1Val = "20251212"; Row2Val = 3;         Row3Val = -26 },
    @{ SheetIdx = 10; Name = "z60";   OldCol = 16; NewCol = 17; Row1Type = "n"; Row1Val = "20251212"; Row2Val = 10;        Row3Val = -82 },
    @{ SheetIdx = 11; Name = "gap";   OldCol = 56; NewCol = 57; Row1Type = "s"; Row1Val = "20251212"; Row2Val = 100;       Row3Val = 84 },
    @{ SheetIdx = 12; Name = "std";   OldCol = 37; NewCol = 38; Row1Type = "s"; Row1Val = "20251212"; Row2Val = 4.55;      Row3Val = -5.74 },
    @{ SheetIdx = 13; Name = "quant"; OldCol = 16; NewCol = 17; Row1Type = "s"; Row1Val = "20251212"; Row2Val = 64;        Row3Val = 99 }
)

foreach ($entry in $entries) {
    $ws = $wb.Worksheets.Item($entry.SheetIdx)

    $oldHeader = $ws.Cells.Item(1, $entry.OldCol)
    $newHeader = $ws.Cells.Item(1, $entry.NewCol)
    $oldRow2   = $ws.Cells.Item(2, $entry.OldCol)
    $newRow2   = $ws.Cells.Item(2, $entry.NewCol)
    $oldRow3   = $ws.Cells.Item(3, $entry.OldCol)
    $newRow3   = $ws.Cells.Item(3, $entry.NewCol)

    # Match column width of the new column to the column being duplicated.
    $ws.Columns.Item($entry.NewCol).ColumnWidth = $ws.Columns.Item($entry.OldCol).ColumnWidth

    # Header cell (row 1): carry over the bold/shaded header formatting.
    $oldHeader.Copy()
    $newHeader.PasteSpecial($xlPasteFormats)
    if ($entry.Row1Type -eq "s") {
        # Source stores the date as text -- write it as a formula that
        # evaluates to the literal string, then flatten to a plain value so
        # the result is a genuine text cell (not a formula, not a
        # number), matching the existing sheet's convention.
        $newHeader.Formula = '="' + $entry.Row1Val + '"'
        $newHeader.Copy()
        $newHeader.PasteSpecial($xlPasteValues)
    } else {
        $newHeader.Value = [double]$entry.Row1Val
    }

    # Row 2 / Row 3 data cells: plain numbers, no special formatting beyond
    # whatever the column above them already carries.
    $oldRow2.Copy()
    $newRow2.PasteSpecial($xlPasteFormats)
    $newRow2.Value = $entry.Row2Val

    $oldRow3.Copy()
    $newRow3.PasteSpecial($xlPasteFormats)
    $newRow3.Value = $entry.Row3Val
}
